# Append a new "Time dBase Arch and Manager" row to the Documentation sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New feature row: Feature # 20, Description "Create a dBase based on Time"
$ws.Range("B23").Value = 20
$ws.Range("C23").Value = "Create a dBase based on Time"

# Match the wrapping style used by the rest of the Description column (C).
$ws.Range("C23").WrapText = $true

# Move the selection to the next empty row, like a user who just finished
# typing the new entry and pressed Enter.
[void]$ws.Range("C24").Select()
